$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-42: update Price (D) and Volume(1h) (E) values ---
# NumberFormat "@" (Text) is applied first on cells whose new Price
# string would otherwise be auto-parsed by Excel as a number, so the
# cell keeps storing the literal text (matches the original inline-string type).
$ws.Range("D2").Value = "60.292.92"
$ws.Range("E2").Value = "  -2.01%  "
$ws.Range("D3").Value = "3.382.30"
$ws.Range("E3").Value = "  -1.80%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.81"
$ws.Range("E5").Value = "  -2.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.40"
$ws.Range("E6").Value = "  -6.13%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.382.40"
$ws.Range("E8").Value = "  -1.82%  "
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.47"
$ws.Range("E10").Value = "  -4.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.124"
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.388"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("D13").Value = "3.955.78"
$ws.Range("E13").Value = "  -1.89%  "
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.94"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").Value = "3.385.35"
$ws.Range("E16").Value = "  -1.72%  "
$ws.Range("E17").Value = "  -3.31%  "
$ws.Range("D18").Value = "60.400.25"
$ws.Range("E18").Value = "  -1.96%  "
$ws.Range("E19").Value = "  -2.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.92"
$ws.Range("E20").Value = "  -2.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.01"
$ws.Range("E21").Value = "  -5.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "386.31"
$ws.Range("E22").Value = "  -0.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.557"
$ws.Range("E23").Value = "  -2.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.23"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000115"
$ws.Range("E26").Value = "  -6.36%  "
$ws.Range("D27").Value = "3.528.85"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("E28").Value = "  -1.85%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.35"
$ws.Range("E30").Value = "  -5.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.89"
$ws.Range("E31").Value = "  -4.54%  "
$ws.Range("E32").Value = "  -2.33%  "
$ws.Range("E33").Value = "  -8.24%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.59"
$ws.Range("E35").Value = "  -1.81%  "
$ws.Range("D36").Value = "3.410.68"
$ws.Range("E36").Value = "  -1.70%  "
$ws.Range("E37").Value = "  -2.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "167.37"
$ws.Range("E38").Value = "  +0.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.50"
$ws.Range("E39").Value = "  -4.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.90"
$ws.Range("E40").Value = "  -7.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0770"
$ws.Range("E41").Value = "  -2.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.00"
$ws.Range("E42").Value = "  +1.50%  "

# --- Rows 43 & 44: FirstDigitalUSD and Mantle swapped places, with updated values ---
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.780"
$ws.Range("E44").Value = "  -1.76%  "

# --- Rows 45-51: update Price (D) and Volume(1h) (E) values ---
$ws.Range("E45").Value = "  -1.93%  "
$ws.Range("E46").Value = "  -1.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.22"
$ws.Range("E47").Value = "  -2.46%  "
$ws.Range("D48").Value = "2.514.54"
$ws.Range("E48").Value = "  -3.76%  "
$ws.Range("E49").Value = "  -4.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.00"
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("E51").Value = "  -3.97%  "
